$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, shifting existing rows 3 and 4 down to 4 and 5.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with a duplicate of row 2's values
# (Study Year = 3, Personal ID = 6001016330261)
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 6001016330261
